# Auto-generated edit script applying the numeric corrections described in the diff.
# For each affected leve row, currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) are updated to the new recalculated values. Cells that no longer hold a value
# after the edit are cleared with ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1409.8636
$ws.Range("I40").Value = 1042.8572
$ws.Range("J40").Value = 1581.1333
$ws.Range("K40").Value = 1042.8572
$ws.Range("L40").Value = 1581.1333
$ws.Range("M40").Value = -867.8571999999999
$ws.Range("N40").Value = -1931.1333

$ws.Range("H125").Value = 125000270
$ws.Range("J125").Value = 250000140
$ws.Range("L125").Value = 2250001260
$ws.Range("N125").Value = -2250006180

$ws.Range("H137").Value = 5322.857
$ws.Range("I137").Value = 6208.1665
$ws.Range("J137").Value = 4142.4443
$ws.Range("K137").Value = 18624.4995
$ws.Range("L137").Value = 12427.3329
$ws.Range("M137").Value = -16074.4995
$ws.Range("N137").Value = -17527.3329

$ws.Range("H138").Value = 1806.8052
$ws.Range("I138").Value = 872.8511
$ws.Range("J138").Value = 3270
$ws.Range("K138").Value = 2618.5533
$ws.Range("L138").Value = 9810
$ws.Range("M138").Value = 2521.4467
$ws.Range("N138").Value = -20090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6262.9795
$ws.Range("I32").Value = 4025.7615
$ws.Range("J32").Value = 28138
$ws.Range("K32").Value = 4025.7615
$ws.Range("L32").Value = 28138
$ws.Range("M32").Value = -3738.7615
$ws.Range("N32").Value = -28712

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H44").Value = 14990
$ws.Range("J44").Value = 14990
$ws.Range("L44").Value = 14990
$ws.Range("N44").Value = -15966

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 243435.64
$ws.Range("I61").Value = 225682.33
$ws.Range("J61").Value = 266328.1
$ws.Range("K61").Value = 225682.33
$ws.Range("L61").Value = 266328.1
$ws.Range("M61").Value = -225470.33
$ws.Range("N61").Value = -266752.1

$ws.Range("H74").Value = 194922.88
$ws.Range("I74").Value = 313495.78
$ws.Range("J74").Value = 43149.56
$ws.Range("K74").Value = 313495.78
$ws.Range("L74").Value = 43149.56
$ws.Range("M74").Value = -312621.78
$ws.Range("N74").Value = -44897.56

$ws.Range("H77").Value = 194922.88
$ws.Range("I77").Value = 313495.78
$ws.Range("J77").Value = 43149.56
$ws.Range("K77").Value = 1567478.9
$ws.Range("L77").Value = 215747.8
$ws.Range("M77").Value = -1563110.9
$ws.Range("N77").Value = -224483.8

$ws.Range("H132").Value = 3347.681
$ws.Range("I132").Value = 3953.2273
$ws.Range("J132").Value = 2814.8
$ws.Range("K132").Value = 11859.6819
$ws.Range("L132").Value = 8444.400000000001
$ws.Range("M132").Value = -9329.6819
$ws.Range("N132").Value = -13504.4

$ws.Range("H136").Value = 243435.64
$ws.Range("I136").Value = 225682.33
$ws.Range("J136").Value = 266328.1
$ws.Range("K136").Value = 677046.99
$ws.Range("L136").Value = 798984.2999999999
$ws.Range("M136").Value = -674496.99
$ws.Range("N136").Value = -804084.2999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3142.761
$ws.Range("I134").Value = 2780.2188
$ws.Range("J134").Value = 3971.4285
$ws.Range("K134").Value = 8340.6564
$ws.Range("L134").Value = 11914.2855
$ws.Range("M134").Value = -5805.6564
$ws.Range("N134").Value = -16984.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2734.0698
$ws.Range("I31").Value = 1389.1562
$ws.Range("J31").Value = 6646.5454
$ws.Range("K31").Value = 1389.1562
$ws.Range("L31").Value = 6646.5454
$ws.Range("M31").Value = -1094.1562
$ws.Range("N31").Value = -7236.5454

$ws.Range("H34").Value = 2734.0698
$ws.Range("I34").Value = 1389.1562
$ws.Range("J34").Value = 6646.5454
$ws.Range("K34").Value = 1389.1562
$ws.Range("L34").Value = 6646.5454
$ws.Range("M34").Value = -1187.1562
$ws.Range("N34").Value = -7050.5454

$ws.Range("H99").Value = 336004
$ws.Range("I99").Value = 336004
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 336004
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -334506
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 246.75
$ws.Range("I107").Value = 246.75
$ws.Range("K107").Value = 246.75
$ws.Range("M107").Value = 1673.25

$ws.Range("H126").Value = 336004
$ws.Range("I126").Value = 336004
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1008012
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1005542
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2187.7837
$ws.Range("I132").Value = 1081.6666
$ws.Range("J132").Value = 4229.846
$ws.Range("K132").Value = 3244.9998
$ws.Range("L132").Value = 12689.538
$ws.Range("M132").Value = -714.9998000000001
$ws.Range("N132").Value = -17749.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 210
$ws.Range("I60").Value = 100
$ws.Range("J60").Value = 375
$ws.Range("K60").Value = 300
$ws.Range("L60").Value = 1125
$ws.Range("M60").Value = -49
$ws.Range("N60").Value = -1627

$ws.Range("H68").Value = 580.9
$ws.Range("I68").Value = 502.25
$ws.Range("J68").Value = 633.3333
$ws.Range("K68").Value = 1506.75
$ws.Range("L68").Value = 1899.9999
$ws.Range("M68").Value = -695.75
$ws.Range("N68").Value = -3521.9999

$ws.Range("H71").Value = 580.9
$ws.Range("I71").Value = 502.25
$ws.Range("J71").Value = 633.3333
$ws.Range("K71").Value = 4520.25
$ws.Range("L71").Value = 5699.9997
$ws.Range("M71").Value = -464.25
$ws.Range("N71").Value = -13811.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 8000
$ws.Range("J28").Value = 8000
$ws.Range("L28").Value = 8000
$ws.Range("N28").Value = -8384

$ws.Range("H93").Value = 26387.285
$ws.Range("J93").Value = 26387.285
$ws.Range("L93").Value = 26387.285
$ws.Range("N93").Value = -30131.285

$ws.Range("H132").Value = 4361.6
$ws.Range("I132").Value = 4483.4443
$ws.Range("J132").Value = 4232.5884
$ws.Range("K132").Value = 13450.3329
$ws.Range("L132").Value = 12697.7652
$ws.Range("M132").Value = -10920.3329
$ws.Range("N132").Value = -17757.7652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2464.2593
$ws.Range("I7").Value = 2241.25
$ws.Range("J7").Value = 2642.6667
$ws.Range("K7").Value = 2241.25
$ws.Range("L7").Value = 2642.6667
$ws.Range("M7").Value = -2129.25
$ws.Range("N7").Value = -2866.6667

$ws.Range("H46").Value = 1655.3334
$ws.Range("I46").Value = 1995
$ws.Range("J46").Value = 976
$ws.Range("K46").Value = 1995
$ws.Range("L46").Value = 976
$ws.Range("M46").Value = -1807
$ws.Range("N46").Value = -1352

$ws.Range("H126").Value = 2464.2593
$ws.Range("I126").Value = 2241.25
$ws.Range("J126").Value = 2642.6667
$ws.Range("K126").Value = 6723.75
$ws.Range("L126").Value = 7928.000100000001
$ws.Range("M126").Value = -4253.75
$ws.Range("N126").Value = -12868.0001

$ws.Range("H132").Value = 4881.808
$ws.Range("I132").Value = 2100.6365
$ws.Range("J132").Value = 9101.518
$ws.Range("K132").Value = 6301.9095
$ws.Range("L132").Value = 27304.554
$ws.Range("M132").Value = -3771.9095
$ws.Range("N132").Value = -32364.554

$ws.Range("H136").Value = 3187.6196
$ws.Range("I136").Value = 1800.6735
$ws.Range("K136").Value = 5402.020500000001
$ws.Range("M136").Value = -2852.020500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 6654.6665
$ws.Range("I26").Value = 2500
$ws.Range("J26").Value = 8732
$ws.Range("K26").Value = 2500
$ws.Range("L26").Value = 8732
$ws.Range("M26").Value = -2207
$ws.Range("N26").Value = -9318

$ws.Range("H29").Value = 70011
$ws.Range("J29").Value = 70011
$ws.Range("L29").Value = 70011
$ws.Range("N29").Value = -70591

$ws.Range("H49").Value = 1670304
$ws.Range("I49").Value = 1670304
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 1670304
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -1670074
$ws.Range("N49").ClearContents()

$ws.Range("H126").Value = 1240.1923
$ws.Range("I126").Value = 1029.7778
$ws.Range("J126").Value = 1713.625
$ws.Range("K126").Value = 3089.3334
$ws.Range("L126").Value = 5140.875
$ws.Range("M126").Value = -619.3334000000004
$ws.Range("N126").Value = -10080.875

$ws.Range("H132").Value = 2394.9546
$ws.Range("I132").Value = 1673.1666
$ws.Range("K132").Value = 5019.4998
$ws.Range("M132").Value = -2489.4998

$ws.Range("H136").Value = 11895715
$ws.Range("I136").Value = 18201264
$ws.Range("J136").Value = 335540.5
$ws.Range("K136").Value = 54603792
$ws.Range("L136").Value = 1006621.5
$ws.Range("M136").Value = -54601242
$ws.Range("N136").Value = -1011721.5
